$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 6
    3  = 5
    4  = 2
    5  = 4
    6  = 4
    7  = 4
    8  = 7
    9  = 8
    10 = 11
    11 = 4
    12 = 9
    13 = 4
    14 = 3
    15 = 11
    16 = 5
    17 = 8
    18 = 8
    19 = 6
    20 = 8
    21 = 5
    22 = 5
    23 = 5
    24 = 12
    25 = 4
    26 = 4
    27 = 7
    28 = 5
    29 = 6
    30 = 7
    31 = 3
    32 = 2
    33 = 5
    34 = 7
    35 = 6
    36 = 12
    37 = 11
    38 = 7
    39 = 10
    40 = 2
    41 = 8
    42 = 3
    43 = 5
    44 = 1
    45 = 1
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
